# Apply cryptos list price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.595.28"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.925.75"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4735"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2915"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06855"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.33%  "

$ws.Range("D12").Value = "1.923.95"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07721"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.333"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6726"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "290.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.08%  "

$ws.Range("D17").Value = "30.624.20"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007667"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.598"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.62%  "

$ws.Range("D22").Value = "2.170.95"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.469"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.539"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.56"
$ws.Range("D26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.136"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1073"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.408"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.204"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05048"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7362"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02068"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.98%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.055"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.71"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4468"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8721"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.903"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.99"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.325"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.380"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1257"
$ws.Range("D49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "48.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +13.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.25"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.95%  "
